# "Adjusted notes for Content Inventory"
# Fill in the NOTES column (E) for the Sign Up, Contact Us and Privacy Policy
# rows, and move the active selection to E16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sign Up row (row 6) - note about the sign up form
$ws.Range("E6").Value = "Sign up form, need to get what information is required"

# Contact Us row (row 11) - note about the contact us form
$ws.Range("E11").Value = "Contact us form"

# Privacy Policy row (row 12) - note about needing a privacy policy
$ws.Range("E12").Value = "Need to obtain a privacy policy"

# Update the selected / active cell shown when the workbook is reopened
$null = $ws.Range("E16").Select()
